# Apply the translation-test workbook fixture changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet from "Sheet1" to "Translations"
$ws.Name = "Translations"

# Rebuild the header row and data rows with the new "Entity Id" column
$ws.Range("A1").Value = "Entity Id"
$ws.Range("B1").Value = "Type"
$ws.Range("C1").Value = "Index"
$ws.Range("D1").Value = "Original"
$ws.Range("E1").Value = "Translation"

$ws.Range("A2").Value = "AAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAA"
$ws.Range("B2").Value = "Title"
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = "Orig"
$ws.Range("E2").Value = "title"

$ws.Range("A3").Value = "AAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAA"
$ws.Range("B3").Value = "ValidationMessage"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "Orig"
$ws.Range("E3").Value = "validation message"

$ws.Range("A4").Value = "AAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAA"
$ws.Range("B4").Value = "Instruction"
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = "Orig"
$ws.Range("E4").Value = "instruction"

$ws.Range("A5").Value = "AAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAA"
$ws.Range("B5").Value = "OptionTitle"
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = "Orig"
$ws.Range("E5").Value = "option"

# Column widths / best-fit sizing for the new layout (values tuned so the
# engine's internal character-width quantization lands on the saved widths)
$ws.Columns.Item(1).ColumnWidth = 42.417
$ws.Columns.Item(2).ColumnWidth = 17.26
$ws.Columns.Item(3).ColumnWidth = 5.084

# Move the active selection to match the saved view state
$ws.Range("E10").Select()
